# Scheduled-runner refresh of per-Leve crafting profit figures (currentAveragePrice /
# LevePrice / LeveProfit columns H:N) across all eight Disciple of the Hand sheets.
# Values below mirror a fresh market-board pull; a few rows also gain/lose an HQ or NQ
# profit cell entirely where HQ/NQ availability flipped.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 634.4568
$ws.Range("I15").Value = 634.4568
$ws.Range("K15").Value = 1903.3704
$ws.Range("M15").Value = -1734.3704
$ws.Range("H17").Value = 1778.0769
$ws.Range("I17").Value = 1469.8
$ws.Range("J17").Value = 1970.75
$ws.Range("K17").Value = 4409.4
$ws.Range("L17").Value = 5912.25
$ws.Range("M17").Value = -4241.4
$ws.Range("N17").Value = -6248.25
$ws.Range("H132").Value = 3265283
$ws.Range("I132").Value = 3591587.8
$ws.Range("K132").Value = 10774763.4
$ws.Range("M132").Value = -10772233.4
$ws.Range("H137").Value = 7197.6616
$ws.Range("I137").Value = 11335.546
$ws.Range("J137").Value = 3296.2285
$ws.Range("K137").Value = 34006.638
$ws.Range("L137").Value = 9888.6855
$ws.Range("M137").Value = -31456.638
$ws.Range("N137").Value = -14988.6855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3108.12
$ws.Range("I61").Value = 2508.8262
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 2508.8262
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -2296.8262
$ws.Range("N61").Value = -10424
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 156619.2
$ws.Range("I74").Value = 190135.97
$ws.Range("K74").Value = 190135.97
$ws.Range("M74").Value = -189261.97
$ws.Range("H77").Value = 156619.2
$ws.Range("I77").Value = 190135.97
$ws.Range("K77").Value = 950679.85
$ws.Range("M77").Value = -946311.85
$ws.Range("H122").Value = 1780.2
$ws.Range("I122").Value = 1227
$ws.Range("J122").Value = 3439.8
$ws.Range("K122").Value = 3681
$ws.Range("L122").Value = 10319.4
$ws.Range("M122").Value = -1231
$ws.Range("N122").Value = -15219.4
$ws.Range("H132").Value = 3212
$ws.Range("J132").Value = 3342.6667
$ws.Range("L132").Value = 10028.0001
$ws.Range("N132").Value = -15088.0001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 3108.12
$ws.Range("I136").Value = 2508.8262
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 7526.4786
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -4976.4786
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4530.75
$ws.Range("I105").Value = 1649.4
$ws.Range("J105").Value = 9333
$ws.Range("K105").Value = 1649.4
$ws.Range("L105").Value = 9333
$ws.Range("M105").Value = 97.59999999999991
$ws.Range("N105").Value = -12827
$ws.Range("H134").Value = 4739.8887
$ws.Range("J134").Value = 4353.25
$ws.Range("L134").Value = 13059.75
$ws.Range("N134").Value = -18129.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4004271.5
$ws.Range("I31").Value = 4765513
$ws.Range("K31").Value = 4765513
$ws.Range("M31").Value = -4765218
$ws.Range("H34").Value = 4004271.5
$ws.Range("I34").Value = 4765513
$ws.Range("K34").Value = 4765513
$ws.Range("M34").Value = -4765311
$ws.Range("H99").Value = 5436.625
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498
$ws.Range("H126").Value = 5436.625
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530
$ws.Range("H132").Value = 41800.535
$ws.Range("I132").Value = 47650.23
$ws.Range("K132").Value = 142950.69
$ws.Range("M132").Value = -140420.69
$ws.Range("H141").Value = 303322.2
$ws.Range("J141").Value = 303322.2
$ws.Range("L141").Value = 303322.2
$ws.Range("N141").Value = -313682.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 2016.4
$ws.Range("I26").Value = 20
$ws.Range("J26").Value = 10002
$ws.Range("K26").Value = 60
$ws.Range("L26").Value = 30006
$ws.Range("M26").Value = 228
$ws.Range("N26").Value = -30582
$ws.Range("H51").Value = 3320.8235
$ws.Range("I51").Value = 99
$ws.Range("J51").Value = 3522.1875
$ws.Range("K51").Value = 297
$ws.Range("L51").Value = 10566.5625
$ws.Range("M51").Value = 163
$ws.Range("N51").Value = -11486.5625
$ws.Range("H57").Value = 6744.1665
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 6744.1665
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 20232.4995
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -21350.4995
$ws.Range("H107").Value = 1087.1428
$ws.Range("I107").Value = 1116.5
$ws.Range("K107").Value = 3349.5
$ws.Range("M107").Value = -1429.5
$ws.Range("H113").Value = 2102.923
$ws.Range("I113").Value = 2149.6365
$ws.Range("J113").Value = 2068.6667
$ws.Range("K113").Value = 6448.9095
$ws.Range("L113").Value = 6206.000100000001
$ws.Range("M113").Value = -4278.9095
$ws.Range("N113").Value = -10546.0001
$ws.Range("H122").Value = 1177.3928
$ws.Range("I122").Value = 818.5
$ws.Range("K122").Value = 7366.5
$ws.Range("M122").Value = -4916.5
$ws.Range("H123").Value = 3354
$ws.Range("I123").Value = 929.5
$ws.Range("K123").Value = 2788.5
$ws.Range("M123").Value = -338.5
$ws.Range("H131").Value = 102936.98
$ws.Range("J131").Value = 1967.3549
$ws.Range("L131").Value = 5902.0647
$ws.Range("N131").Value = -15982.0647

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1023.6667
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H63").Value = 33333.332
$ws.Range("J63").Value = 33333.332
$ws.Range("L63").Value = 33333.332
$ws.Range("N63").Value = -34705.332
$ws.Range("H66").Value = 33333.332
$ws.Range("J66").Value = 33333.332
$ws.Range("L66").Value = 99999.99600000001
$ws.Range("N66").Value = -106863.996
$ws.Range("H80").Value = 4663.609
$ws.Range("I80").Value = 2897.9412
$ws.Range("J80").Value = 9666.333000000001
$ws.Range("K80").Value = 2897.9412
$ws.Range("L80").Value = 9666.333000000001
$ws.Range("M80").Value = -1899.9412
$ws.Range("N80").Value = -11662.333
$ws.Range("H83").Value = 4663.609
$ws.Range("I83").Value = 2897.9412
$ws.Range("J83").Value = 9666.333000000001
$ws.Range("K83").Value = 14489.706
$ws.Range("L83").Value = 48331.665
$ws.Range("M83").Value = -9497.706000000002
$ws.Range("N83").Value = -58315.665
$ws.Range("H122").Value = 3641.1365
$ws.Range("I122").Value = 3676
$ws.Range("K122").Value = 11028
$ws.Range("M122").Value = -8578
$ws.Range("H136").Value = 8463.826999999999
$ws.Range("J136").Value = 8463.826999999999
$ws.Range("L136").Value = 25391.481
$ws.Range("N136").Value = -30491.481

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 11980.765
$ws.Range("I82").Value = 9431.333000000001
$ws.Range("K82").Value = 9431.333000000001
$ws.Range("M82").Value = -9070.333000000001
$ws.Range("H85").Value = 11980.765
$ws.Range("I85").Value = 9431.333000000001
$ws.Range("K85").Value = 9431.333000000001
$ws.Range("M85").Value = -8183.333000000001
$ws.Range("H122").Value = 12143.72
$ws.Range("I122").Value = 11687.45
$ws.Range("K122").Value = 35062.35000000001
$ws.Range("M122").Value = -32612.35000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 19222.295
$ws.Range("I132").Value = 24317.883
$ws.Range("J132").Value = 1897.3
$ws.Range("K132").Value = 72953.649
$ws.Range("L132").Value = 5691.9
$ws.Range("M132").Value = -70423.649
$ws.Range("N132").Value = -10751.9
